$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) date serial from 45184 to 45185 for rows 2-10
$ws.Range("C2:C10").Value = 45185
